# Adds rows 98-107 to Sheet1 with Defense-gov Explore Feed article data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('A98').Value = 'Defense-gov Explore Feed'
$ws.Range('B98').Value = 'Off-Camera, On-Record Press Briefing'
$ws.Range('C98').Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3497547/off-camera-on-record-press-briefing/'
$ws.Cells.Item(98, 4).NumberFormat = '@'
$ws.Cells.Item(98, 4).Value = '2023-08-18'
$ws.Cells.Item(98, 4).ClearFormats()
$ws.Range('E98').Value = 'nt Task Force 5-0 Dual Status Commander Brig. Gen. Stephen F. Logan will hold a virtual briefing with updates on the DOD''s response to the Maui wildfires.'

$ws.Range('A99').Value = 'Defense-gov Explore Feed'
$ws.Range('B99').Value = 'Media Roundtable: Military Service Academies On-Site Installation Evaluations (MSA OSIE) Today'
$ws.Range('C99').Value = 'https://www.defense.gov/News/Advisories/Advisory/Article/3495874/media-roundtable-military-service-academies-on-site-installation-evaluations-ms/'
$ws.Cells.Item(99, 4).NumberFormat = '@'
$ws.Cells.Item(99, 4).Value = '2023-08-17'
$ws.Cells.Item(99, 4).ClearFormats()
$ws.Range('E99').Value = ' Elizabeth Foster, Executive Director of the Office of Force Resiliency for the Under Secretary of Defense for Personnel and Readiness and Response Office, and Dr. Andra Tharp, the senior prevention advisor for the department''s Office of Force Resiliency and SAPRO will hold an off-camera, on-the-record media roundtable.'

$ws.Range('A100').Value = 'Defense-gov Explore Feed'
$ws.Range('B100').Value = 'U.S. Department of Defense and Japan Ministry of Defense Press Release on the Commencement of Glide Phase Interceptor Cooperative Development'
$ws.Range('C100').Value = 'https://www.defense.gov/News/Releases/Release/Article/3498431/us-department-of-defense-and-japan-ministry-of-defense-press-release-on-the-com/'
$ws.Cells.Item(100, 4).NumberFormat = '@'
$ws.Cells.Item(100, 4).Value = '2023-08-18'
$ws.Cells.Item(100, 4).ClearFormats()
$ws.Range('E100').Value = 'ing the January 2023 Security Consultative Committee (“2+2”) meeting, the United States and Japan concurred with beginning discussions on potential joint development of a future interceptor.'

$ws.Range('A101').Value = 'Defense-gov Explore Feed'
$ws.Range('B101').Value = 'Statement from Secretary of Defense Lloyd J. Austin III on United States Trilateral Summit with Japan and the Republic of Korea'
$ws.Range('C101').Value = 'https://www.defense.gov/News/Releases/Release/Article/3498219/statement-from-secretary-of-defense-lloyd-j-austin-iii-on-united-states-trilate/'
$ws.Cells.Item(101, 4).NumberFormat = '@'
$ws.Cells.Item(101, 4).Value = '2023-08-18'
$ws.Cells.Item(101, 4).ClearFormats()
$ws.Range('E101').Value = 'tatement from Secretary of Defense Lloyd J. Austin III on the United States trilateral summit with Japan and South Korea.'

$ws.Range('A102').Value = 'Defense-gov Explore Feed'
$ws.Range('B102').Value = 'Honorable Mara Karlin Concludes Visit to Djibouti'
$ws.Range('C102').Value = 'https://www.defense.gov/News/Releases/Release/Article/3497932/honorable-mara-karlin-concludes-visit-to-djibouti/'
$ws.Cells.Item(102, 4).NumberFormat = '@'
$ws.Cells.Item(102, 4).Value = '2023-08-18'
$ws.Cells.Item(102, 4).ClearFormats()
$ws.Range('E102').Value = 'a Karlin, performing the duties of deputy undersecretary of defense for policy, met with U.S. Africa Command and Combined Joint Task Force Horn of Africa leaders during a visit to Djibouti.'

$ws.Range('A103').Value = 'Defense-gov Explore Feed'
$ws.Range('B103').Value = 'Navy to Christen Guided-Missile Destroyer Ted Stevens (DDG 128)'
$ws.Range('C103').Value = 'https://www.defense.gov/News/Releases/Release/Article/3497493/navy-to-christen-guided-missile-destroyer-ted-stevens-ddg-128/'
$ws.Cells.Item(103, 4).NumberFormat = '@'
$ws.Cells.Item(103, 4).Value = '2023-08-18'
$ws.Cells.Item(103, 4).ClearFormats()
$ws.Range('E103').Value = ' Navy will christen the future USS Ted Stevens during a 9:00 a.m. CDT ceremony on Saturday, Aug. 19, in Pascagoula, Mississippi.'

$ws.Range('A104').Value = 'Defense-gov Explore Feed'
$ws.Range('B104').Value = 'Department of Defense Releases Actions to Transform Climate and Enhance Prevention of Harmful Behaviors at Military Service Academies'
$ws.Range('C104').Value = 'https://www.defense.gov/News/Releases/Release/Article/3496394/department-of-defense-releases-actions-to-transform-climate-and-enhance-prevent/'
$ws.Cells.Item(104, 4).NumberFormat = '@'
$ws.Cells.Item(104, 4).Value = '2023-08-17'
$ws.Cells.Item(104, 4).ClearFormats()
$ws.Range('E104').Value = 'retary of Defense Lloyd J. Austin III directed multiple actions to transform climate and enhance prevention of harmful behaviors at the Military Service Academies.'

$ws.Range('A105').Value = 'Defense-gov Explore Feed'
$ws.Range('B105').Value = 'Department of Defense-Sponsored Cyber Internship Offers Knowledge, Inspiration for College Students'
$ws.Range('C105').Value = 'https://www.defense.gov/News/Releases/Release/Article/3496391/department-of-defense-sponsored-cyber-internship-offers-knowledge-inspiration-f/'
$ws.Cells.Item(105, 4).NumberFormat = '@'
$ws.Cells.Item(105, 4).Value = '2023-08-17'
$ws.Cells.Item(105, 4).ClearFormats()
$ws.Range('E105').Value = ' Office of the Under Secretary of Defense for Research and Engineering-sponsored Cyber-Spectrum internship program, MAVEN, recently celebrated its second internship graduation.'

$ws.Range('A106').Value = 'Defense-gov Explore Feed'
$ws.Range('B106').Value = 'DOD Unveils Collaborative Biodefense Reforms in Posture Review'
$ws.Range('C106').Value = 'https://www.defense.gov/News/Releases/Release/Article/3495836/dod-unveils-collaborative-biodefense-reforms-in-posture-review/'
$ws.Cells.Item(106, 4).NumberFormat = '@'
$ws.Cells.Item(106, 4).Value = '2023-08-17'
$ws.Cells.Item(106, 4).ClearFormats()
$ws.Range('E106').Value = ' DOD released the Biodefense Posture Review, outlining reforms aimed to posture the department to fight and win in the face of future biothreats.'

$ws.Range('A107').Value = 'Defense-gov Explore Feed'
$ws.Range('B107').Value = 'Deputy Secretary of Defense Kathleen Hicks Statement on Blended-Wing-Body Aircraft Prototype'
$ws.Range('C107').Value = 'https://www.defense.gov/News/Releases/Release/Article/3496004/deputy-secretary-of-defense-kathleen-hicks-statement-on-blended-wing-body-aircr/'
$ws.Cells.Item(107, 4).NumberFormat = '@'
$ws.Cells.Item(107, 4).Value = '2023-08-17'
$ws.Cells.Item(107, 4).ClearFormats()
$ws.Range('E107').Value = ' Defense Department awarded a contract to JetZero for the next phase of a blended-wing body aircraft prototype.'
